$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# --- Remove the existing _GoBack bookmark first (it's about to be
#     recreated at the point of the most recent text edit below); doing
#     this before creating the new one avoids any ambiguity between the
#     two same-named bookmarks.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# --- Change 1: split the "...weren't fixed." run, inserting a _GoBack
#     bookmark (Word's auto-tracked last-edit marker) between "fixe" and
#     "d.", while leaving the following " My damages..." run intact.
$find1 = $d.Content.Find
$find1.Execute("For breach of the warranty of habitability")
$para1 = $find1.Parent.Paragraphs(1)
$range1 = $para1.Range
$textEnd1 = $range1.End - 1   # exclude the paragraph mark
$body1 = $d.Range($range1.Start, $textEnd1)
$body1.Delete()

$xml1 = "<w:p xmlns:w='$wNs'>" +
  "<w:r><w:t>For breach of the warranty of habitability, I am entitled to damages because my rental unit was worth less due to the condition of my home and the problems that weren" +
  [char]0x2019 + "t fixe</w:t></w:r>" +
  "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
  "<w:r><w:t>d.</w:t></w:r>" +
  "<w:r><w:t xml:space='preserve'> My damages are based on the percentage reduction in the value of my home caused by the landlord" +
  [char]0x2019 + "s failure to make repairs.</w:t></w:r>" +
  "</w:p>"
$insertPoint1 = $d.Range($range1.Start, $range1.Start)
$insertPoint1.InsertXML($xml1)

# --- Change 2: the paragraph that used to hold the _GoBack bookmark
#     (now removed above) gets the relocation/hotel request text.
$find2 = $d.Content.Find
$find2.Execute("{%p if complaint_ask_for_relocation %}")
$para2 = $find2.Parent.Paragraphs(1)
$nextPara = $para2.Next()
$range2 = $nextPara.Range
$insertPoint2 = $d.Range($range2.Start, $range2.Start)
$xml2 = "<w:p xmlns:w='$wNs'><w:r><w:t>If Tenant is required to move out in order to make repairs, order the Landlord to provide alternative housing or to arrange and pay for a hotel that is comparable in size, amenities and location until such time as the repairs have been completed.</w:t></w:r></w:p>"
$insertPoint2.InsertXML($xml2)
